# Add figure paths for the publications that were missing them, and fix the
# casing of the gauge-fixing figure's file extension, per the commit
# "add more figures for selected publications".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = "/pictures/figures/fig1_gaugefixing_posfai_et_al.PNG"
$ws.Range("F6").Value  = "/pictures/figures/fig1_mavenn.png"
$ws.Range("F8").Value  = "/pictures/figures/fig1_mpra_review.png"
$ws.Range("F9").Value  = "/pictures/figures/fig1_density_estimation.png"
$ws.Range("F10").Value = "/pictures/figures/fig1_tite_seq.png"
$ws.Range("F11").Value = "/pictures/figures/fig1_equitability_kinney_2014.png"
$ws.Range("F2").Value  = "/pictures/figures/fig1_rousseau.png"

# Update the active selection to match the edited workbook (F3).
$ws.Range("F3").Select()
